$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.263.04"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.232.30"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'243.36"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").Value = "'74.55"
$ws.Range("E7").Value = "  -2.99%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Value = "'42.77"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("D11").Value = "'0.0962"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "'7.00"
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "2.570.43"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "'14.34"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Value = "'0.840"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").Value = "2.235.50"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "42.163.05"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'0.0000106"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "'6.23"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'73.09"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").Value = "'11.36"
$ws.Range("D23").Value = "'230.88"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E24").Value = "  -7.34%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "'11.45"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "'166.87"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "'20.61"
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").Value = "'30.24"
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  -8.70%  "
$ws.Range("D37").Value = "'4.39"
$ws.Range("E37").Value = "  -8.79%  "
$ws.Range("D38").Value = "'0.0305"
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("D39").Value = "'13.28"
$ws.Range("E39").Value = "  -6.82%  "
$ws.Range("D40").Value = "'2.14"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").Value = "'5.70"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "'65.09"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").Value = "'8.73"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "'104.60"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").Value = "'0.100"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "'2.36"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "2.444.53"
$ws.Range("E51").Value = "  -0.92%  "
